# Elimina EC anteriores y se agregan nuevos, se modifica base de datos
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C16").Value = "1047402961"
$ws.Range("D16").Value = "JONATHAN MEZA BUSTAMANTE"
$ws.Range("E16").Value = "2304"
$ws.Range("F16").Value = 46400
$ws.Range("G16").Value = 1160000

$ws.Range("C17").Value = "1047456993"
$ws.Range("D17").Value = "ESTEBAN DE JESUS AHUMEDO BURGOS"
$ws.Range("E17").Value = "2204"
$ws.Range("F17").Value = 8000
$ws.Range("G17").Value = 908526

$ws.Range("C18").Value = "73162704"
$ws.Range("D18").Value = "WILSON FRIAS ALCALA"
$ws.Range("E18").Value = "2112"
$ws.Range("F18").Value = 36341
$ws.Range("G18").Value = 908526

$ws.Range("C19").Value = "73119747"
$ws.Range("D19").Value = "PEDRO CLAVER CARABALLO OROZCO"
$ws.Range("E19").Value = "2304"
$ws.Range("F19").Value = 46400
$ws.Range("G19").Value = 1160000

$ws.Range("C20").Value = "9296373"
$ws.Range("D20").Value = "RALLPH ANTONIO LLOREDA DORIA"
$ws.Range("E20").Value = "2304"
$ws.Range("F20").Value = 46400
$ws.Range("G20").Value = 1160000

$ws.Range("C21").Value = "73574969"
$ws.Range("D21").Value = "JOHN ALEXANDER PEREZ DE CASTRO"
$ws.Range("E21").Value = "2302"
$ws.Range("F21").Value = 46400
$ws.Range("G21").Value = 1160000

$ws.Range("C22").Value = "73574969"
$ws.Range("D22").Value = "JOHN ALEXANDER PEREZ DE CASTRO"
$ws.Range("E22").Value = "2209"
$ws.Range("F22").Value = 36341
$ws.Range("G22").Value = 1160000

$ws.Range("C23").Value = "73574969"
$ws.Range("D23").Value = "JOHN ALEXANDER PEREZ DE CASTRO"
$ws.Range("E23").Value = "2208"
$ws.Range("F23").Value = 36341
$ws.Range("G23").Value = 1160000

$ws.Range("C24").Value = "73574969"
$ws.Range("D24").Value = "JOHN ALEXANDER PEREZ DE CASTRO"
$ws.Range("E24").Value = "2207"
$ws.Range("F24").Value = 36341
$ws.Range("G24").Value = 1160000

$ws.Range("C25").Value = "73574969"
$ws.Range("D25").Value = "JOHN ALEXANDER PEREZ DE CASTRO"
$ws.Range("E25").Value = "2206"
$ws.Range("F25").Value = 36341
$ws.Range("G25").Value = 1160000

$ws.Range("C26").Value = "73574969"
$ws.Range("D26").Value = "JOHN ALEXANDER PEREZ DE CASTRO"
$ws.Range("E26").Value = "2205"
$ws.Range("F26").Value = 36341
$ws.Range("G26").Value = 1160000

$ws.Range("C27").Value = "73574969"
$ws.Range("D27").Value = "JOHN ALEXANDER PEREZ DE CASTRO"
$ws.Range("E27").Value = "2204"
$ws.Range("F27").Value = 36341
$ws.Range("G27").Value = 1160000

$ws.Range("C28").Value = "73574969"
$ws.Range("D28").Value = "JOHN ALEXANDER PEREZ DE CASTRO"
$ws.Range("E28").Value = "2203"
$ws.Range("F28").Value = 36341
$ws.Range("G28").Value = 1160000

$ws.Range("C29").Value = "73574969"
$ws.Range("D29").Value = "JOHN ALEXANDER PEREZ DE CASTRO"
$ws.Range("E29").Value = "2202"
$ws.Range("F29").Value = 36341
$ws.Range("G29").Value = 1160000

$ws.Range("C30").Value = "73574969"
$ws.Range("D30").Value = "JOHN ALEXANDER PEREZ DE CASTRO"
$ws.Range("E30").Value = "2201"
$ws.Range("F30").Value = 36341
$ws.Range("G30").Value = 1160000

$ws.Range("C31").Value = "73167930"
$ws.Range("D31").Value = "OMAR DE JESUS TORRES CASTILLO"
$ws.Range("E31").Value = "2304"
$ws.Range("F31").Value = 46400
$ws.Range("G31").Value = 1160000
